$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156 (pushes existing rows 156-302 down to 157-303)
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new weekly record
$ws.Cells.Item(156, 1).Value = 8
$ws.Cells.Item(156, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44790
$ws.Cells.Item(156, 5).Value = 4
$ws.Cells.Item(156, 6).Value = 100112012
$ws.Cells.Item(156, 7).Value = "Espinaca"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 2400
$ws.Cells.Item(156, 11).Value = 500
$ws.Cells.Item(156, 12).Value = 550
$ws.Cells.Item(156, 13).Value = 525
$ws.Cells.Item(156, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(156, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value = 1050
$ws.Cells.Item(156, 17).Value = 0.5
$ws.Cells.Item(156, 18).Value = "Hortaliza"
